$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,4).Value = 44505
$ws.Cells.Item(2,12).Value = "Primera"
$ws.Cells.Item(2,13).Value = 100
$ws.Cells.Item(2,14).Value = 2200
$ws.Cells.Item(2,15).Value = 2200
$ws.Cells.Item(2,16).Value = 2200
$ws.Cells.Item(2,17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(2,18).Value = "Provincia de Limarí"
$ws.Cells.Item(2,19).Value = 2200
$ws.Cells.Item(2,20).Value = 1

# Row 3
$ws.Cells.Item(3,4).Value = 44505
$ws.Cells.Item(3,12).Value = "Segunda"
$ws.Cells.Item(3,13).Value = 100
$ws.Cells.Item(3,14).Value = 1800
$ws.Cells.Item(3,15).Value = 1800
$ws.Cells.Item(3,16).Value = 1800
$ws.Cells.Item(3,17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(3,18).Value = "Provincia de Limarí"
$ws.Cells.Item(3,19).Value = 1800
$ws.Cells.Item(3,20).Value = 1

# Row 4
$ws.Cells.Item(4,4).Value = 44483
$ws.Cells.Item(4,12).Value = "Primera"
$ws.Cells.Item(4,13).Value = 50
$ws.Cells.Item(4,14).Value = 2600
$ws.Cells.Item(4,15).Value = 2600
$ws.Cells.Item(4,16).Value = 2600
$ws.Cells.Item(4,17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(4,18).Value = "Provincia de Limarí"
$ws.Cells.Item(4,19).Value = 2600
$ws.Cells.Item(4,20).Value = 1

# Row 5
$ws.Cells.Item(5,4).Value = 44483
$ws.Cells.Item(5,12).Value = "Segunda"
$ws.Cells.Item(5,13).Value = 50
$ws.Cells.Item(5,14).Value = 2400
$ws.Cells.Item(5,15).Value = 2400
$ws.Cells.Item(5,16).Value = 2400
$ws.Cells.Item(5,17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(5,18).Value = "Provincia de Limarí"
$ws.Cells.Item(5,19).Value = 2400
$ws.Cells.Item(5,20).Value = 1

# Row 6
$ws.Cells.Item(6,4).Value = 44484
$ws.Cells.Item(6,12).Value = "Primera"
$ws.Cells.Item(6,13).Value = 100
$ws.Cells.Item(6,14).Value = 25000
$ws.Cells.Item(6,15).Value = 26000
$ws.Cells.Item(6,16).Value = 25500
$ws.Cells.Item(6,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(6,18).Value = "Provincia de Limarí"
$ws.Cells.Item(6,19).Value = 2550
$ws.Cells.Item(6,20).Value = 10

# Row 7
$ws.Cells.Item(7,4).Value = 44488
$ws.Cells.Item(7,12).Value = "Primera"
$ws.Cells.Item(7,13).Value = 50
$ws.Cells.Item(7,14).Value = 25000
$ws.Cells.Item(7,15).Value = 26000
$ws.Cells.Item(7,16).Value = 25600
$ws.Cells.Item(7,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(7,18).Value = "Provincia de Limarí"
$ws.Cells.Item(7,19).Value = 2560
$ws.Cells.Item(7,20).Value = 10

# Row 8
$ws.Cells.Item(8,4).Value = 44469
$ws.Cells.Item(8,12).Value = "Primera"
$ws.Cells.Item(8,13).Value = 100
$ws.Cells.Item(8,14).Value = 28000
$ws.Cells.Item(8,15).Value = 29000
$ws.Cells.Item(8,16).Value = 28500
$ws.Cells.Item(8,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(8,18).Value = "Provincia de Limarí"
$ws.Cells.Item(8,19).Value = 2850
$ws.Cells.Item(8,20).Value = 10

# Row 9
$ws.Cells.Item(9,4).Value = 44491
$ws.Cells.Item(9,12).Value = "Primera"
$ws.Cells.Item(9,13).Value = 150
$ws.Cells.Item(9,14).Value = 25000
$ws.Cells.Item(9,15).Value = 26000
$ws.Cells.Item(9,16).Value = 25467
$ws.Cells.Item(9,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(9,18).Value = "Provincia de Limarí"
$ws.Cells.Item(9,19).Value = 2547
$ws.Cells.Item(9,20).Value = 10

# Row 10
$ws.Cells.Item(10,4).Value = 44511
$ws.Cells.Item(10,12).Value = "Primera"
$ws.Cells.Item(10,13).Value = 80
$ws.Cells.Item(10,14).Value = 25000
$ws.Cells.Item(10,15).Value = 26000
$ws.Cells.Item(10,16).Value = 25375
$ws.Cells.Item(10,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(10,18).Value = "Provincia de Limarí"
$ws.Cells.Item(10,19).Value = 2538
$ws.Cells.Item(10,20).Value = 10

# Row 11
$ws.Cells.Item(11,4).Value = 44461
$ws.Cells.Item(11,12).Value = "Primera"
$ws.Cells.Item(11,13).Value = 100
$ws.Cells.Item(11,14).Value = 29000
$ws.Cells.Item(11,15).Value = 30000
$ws.Cells.Item(11,16).Value = 29500
$ws.Cells.Item(11,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(11,18).Value = "Provincia de Limarí"
$ws.Cells.Item(11,19).Value = 2950
$ws.Cells.Item(11,20).Value = 10

# Row 12
$ws.Cells.Item(12,4).Value = 44446
$ws.Cells.Item(12,12).Value = "Primera"
$ws.Cells.Item(12,13).Value = 100
$ws.Cells.Item(12,14).Value = 3200
$ws.Cells.Item(12,15).Value = 3300
$ws.Cells.Item(12,16).Value = 3250
$ws.Cells.Item(12,17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(12,18).Value = "Provincia del Elquí"
$ws.Cells.Item(12,19).Value = 3250
$ws.Cells.Item(12,20).Value = 1

# Row 13
$ws.Cells.Item(13,4).Value = 44495
$ws.Cells.Item(13,12).Value = "Primera"
$ws.Cells.Item(13,13).Value = 100
$ws.Cells.Item(13,14).Value = 26000
$ws.Cells.Item(13,15).Value = 27000
$ws.Cells.Item(13,16).Value = 26500
$ws.Cells.Item(13,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(13,18).Value = "Provincia de Limarí"
$ws.Cells.Item(13,19).Value = 2650
$ws.Cells.Item(13,20).Value = 10

# Row 14
$ws.Cells.Item(14,4).Value = 44462
$ws.Cells.Item(14,12).Value = "Primera"
$ws.Cells.Item(14,13).Value = 100
$ws.Cells.Item(14,14).Value = 2900
$ws.Cells.Item(14,15).Value = 3000
$ws.Cells.Item(14,16).Value = 2950
$ws.Cells.Item(14,17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(14,18).Value = "Provincia de Limarí"
$ws.Cells.Item(14,19).Value = 2950
$ws.Cells.Item(14,20).Value = 1

# Row 15
$ws.Cells.Item(15,4).Value = 44462
$ws.Cells.Item(15,12).Value = "Segunda"
$ws.Cells.Item(15,13).Value = 50
$ws.Cells.Item(15,14).Value = 2600
$ws.Cells.Item(15,15).Value = 2600
$ws.Cells.Item(15,16).Value = 2600
$ws.Cells.Item(15,17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(15,18).Value = "Provincia de Limarí"
$ws.Cells.Item(15,19).Value = 2600
$ws.Cells.Item(15,20).Value = 1

# Row 16
$ws.Cells.Item(16,4).Value = 44530
$ws.Cells.Item(16,12).Value = "Primera"
$ws.Cells.Item(16,13).Value = 100
$ws.Cells.Item(16,14).Value = 2000
$ws.Cells.Item(16,15).Value = 2100
$ws.Cells.Item(16,16).Value = 2050
$ws.Cells.Item(16,17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(16,18).Value = "Provincia de Limarí"
$ws.Cells.Item(16,19).Value = 2050
$ws.Cells.Item(16,20).Value = 1

# Row 17
$ws.Cells.Item(17,4).Value = 44475
$ws.Cells.Item(17,12).Value = "Primera"
$ws.Cells.Item(17,13).Value = 100
$ws.Cells.Item(17,14).Value = 29000
$ws.Cells.Item(17,15).Value = 30000
$ws.Cells.Item(17,16).Value = 29500
$ws.Cells.Item(17,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(17,18).Value = "Provincia de Limarí"
$ws.Cells.Item(17,19).Value = 2950
$ws.Cells.Item(17,20).Value = 10

# Row 18
$ws.Cells.Item(18,4).Value = 44160
$ws.Cells.Item(18,12).Value = "Primera"
$ws.Cells.Item(18,13).Value = 100
$ws.Cells.Item(18,14).Value = 17000
$ws.Cells.Item(18,15).Value = 18000
$ws.Cells.Item(18,16).Value = 17500
$ws.Cells.Item(18,17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(18,18).Value = "Provincia de Limarí"
$ws.Cells.Item(18,19).Value = 2188
$ws.Cells.Item(18,20).Value = 8

# Row 19
$ws.Cells.Item(19,4).Value = 44160
$ws.Cells.Item(19,12).Value = "Segunda"
$ws.Cells.Item(19,13).Value = 50
$ws.Cells.Item(19,14).Value = 15000
$ws.Cells.Item(19,15).Value = 15000
$ws.Cells.Item(19,16).Value = 15000
$ws.Cells.Item(19,17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(19,18).Value = "Provincia de Limarí"
$ws.Cells.Item(19,19).Value = 1875
$ws.Cells.Item(19,20).Value = 8

# Row 20
$ws.Cells.Item(20,4).Value = 44516
$ws.Cells.Item(20,12).Value = "Primera"
$ws.Cells.Item(20,13).Value = 100
$ws.Cells.Item(20,14).Value = 1900
$ws.Cells.Item(20,15).Value = 2000
$ws.Cells.Item(20,16).Value = 1950
$ws.Cells.Item(20,17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(20,18).Value = "Provincia de Limarí"
$ws.Cells.Item(20,19).Value = 1950
$ws.Cells.Item(20,20).Value = 1

# Row 21
$ws.Cells.Item(21,4).Value = 44516
$ws.Cells.Item(21,12).Value = "Segunda"
$ws.Cells.Item(21,13).Value = 50
$ws.Cells.Item(21,14).Value = 1700
$ws.Cells.Item(21,15).Value = 1700
$ws.Cells.Item(21,16).Value = 1700
$ws.Cells.Item(21,17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(21,18).Value = "Provincia de Limarí"
$ws.Cells.Item(21,19).Value = 1700
$ws.Cells.Item(21,20).Value = 1

# Row 22
$ws.Cells.Item(22,4).Value = 44454
$ws.Cells.Item(22,12).Value = "Primera"
$ws.Cells.Item(22,13).Value = 100
$ws.Cells.Item(22,14).Value = 30000
$ws.Cells.Item(22,15).Value = 31000
$ws.Cells.Item(22,16).Value = 30500
$ws.Cells.Item(22,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(22,18).Value = "Provincia de Limarí"
$ws.Cells.Item(22,19).Value = 3050
$ws.Cells.Item(22,20).Value = 10

# Row 23
$ws.Cells.Item(23,4).Value = 44517
$ws.Cells.Item(23,12).Value = "Primera"
$ws.Cells.Item(23,13).Value = 100
$ws.Cells.Item(23,14).Value = 25000
$ws.Cells.Item(23,15).Value = 27000
$ws.Cells.Item(23,16).Value = 26000
$ws.Cells.Item(23,17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(23,18).Value = "Provincia de Limarí"
$ws.Cells.Item(23,19).Value = 2600
$ws.Cells.Item(23,20).Value = 10

# Row 24
$ws.Cells.Item(24,4).Value = 44467
$ws.Cells.Item(24,12).Value = "Primera"
$ws.Cells.Item(24,13).Value = 100
$ws.Cells.Item(24,14).Value = 2700
$ws.Cells.Item(24,15).Value = 2800
$ws.Cells.Item(24,16).Value = 2750
$ws.Cells.Item(24,17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(24,18).Value = "Provincia de Limarí"
$ws.Cells.Item(24,19).Value = 2750
$ws.Cells.Item(24,20).Value = 1

# Row 25
$ws.Cells.Item(25,4).Value = 44467
$ws.Cells.Item(25,12).Value = "Segunda"
$ws.Cells.Item(25,13).Value = 50
$ws.Cells.Item(25,14).Value = 2500
$ws.Cells.Item(25,15).Value = 2500
$ws.Cells.Item(25,16).Value = 2500
$ws.Cells.Item(25,17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(25,18).Value = "Provincia de Limarí"
$ws.Cells.Item(25,19).Value = 2500
$ws.Cells.Item(25,20).Value = 1
